# Update: Added Notifications, reworked all events
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 4 (event updated from "Ein Tag vor Berkos Bday" to "Mamas Geburtstag") ---
# All of these columns are stored as text in the workbook (IDs exceed double
# precision and values like "08" need their leading zero preserved), so force
# the Text number format before assigning the values.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "7272754151794020845"

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "08"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "11"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "2022"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "Mamas Geburtstag"

$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "10"

# --- Add new row 5 (new event "Adrianas Geburtstag") ---
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "8252115886235587053"

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "13"

$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "11"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "2022"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "Adrianas Geburtstag"

$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "10"
